$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D23").Value = "안녕하세요! 질문이 있습니다. 노트북 쓰다보면 패키지를 인스톨해서 쓰는 경우가 있는데, 매번 패키지를 인스톨하는걸 피할 수 있는 방법이 있나요?"
$ws.Range("E23").Value = "https://theonly1.tistory.com/2702"

$ws.Range("D32").Value = "카테고리형 변수 변환 : Gumbel Softmax"
$ws.Range("E32").Value = "https://dodonam.tistory.com/296"

$ws.Range("D37").Value = "[Paper Review] Geometic Graph Convolutional Networks"
$ws.Range("E37").Value = "http://dsba.korea.ac.kr/seminar/?uid=1439&mod=document&pageid=1"

$ws.Range("D39").Value = "Facial Landmarks for Face Recognition with Dlib"
$ws.Range("E39").Value = "https://a292run.tistory.com/entry/Facial-Landmarks-for-Face-Recognition-with-Dlib-1"
